$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132, shifting the existing rows 132:151 down to 133:152.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new weekly price-report entry.
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(132, 3).Value = "Ñuble"
$ws.Cells.Item(132, 4).Value = 44474
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112043
$ws.Cells.Item(132, 7).Value = "Pepino ensalada"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 120
$ws.Cells.Item(132, 11).Value = 18000
$ws.Cells.Item(132, 12).Value = 19000
$ws.Cells.Item(132, 13).Value = 18500
$ws.Cells.Item(132, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(132, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(132, 16).Value = 308
$ws.Cells.Item(132, 17).Value = 60
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# Match the date-formatted number format used by the rest of column D.
$ws.Cells.Item(132, 4).NumberFormat = $ws.Cells.Item(133, 4).NumberFormat()
